# Update countries & provincias Spain
# - Re-order a block of low-count countries (Botsuana..Timor Oriental) so the
#   list reads: Botsuana, Republica de Africa Central, Somalia, Liberia,
#   Belice, Islas Virgenes Britanicas, Anguila, Timor Oriental
# - Refresh a handful of per-country case counts (Estados Unidos, Canada,
#   Tunez and the last row of the re-ordered block)
# - Bump the "datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the Botsuana..Timor Oriental block (rows 197-204) -----------
$ws.Cells.Item(197, 1).Value = "Botsuana"
$ws.Cells.Item(198, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(199, 1).Value = "Somalia"
$ws.Cells.Item(200, 1).Value = "Liberia"
$ws.Cells.Item(201, 1).Value = "Belice"
$ws.Cells.Item(202, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(203, 1).Value = "Anguila"
$ws.Cells.Item(204, 1).Value = "Timor Oriental"

# The country that now lands on row 202 (Islas Virgenes Britanicas) gets
# refreshed case numbers (was 2/0/0/2/0/0/0 -> now 3/1/0/3/0/0/0)
$ws.Cells.Item(202, 2).Value = 3   # Casos totales
$ws.Cells.Item(202, 3).Value = 1   # Nuevos casos
$ws.Cells.Item(202, 4).Value = 0   # Casos activos
$ws.Cells.Item(202, 5).Value = 3   # Recuperados
$ws.Cells.Item(202, 6).Value = 0   # Casos criticos
$ws.Cells.Item(202, 7).Value = 0   # Muertes hoy
$ws.Cells.Item(202, 8).Value = 0   # Muertes

# --- Refresh per-country totals --------------------------------------------
# Row 4: Estados Unidos
$ws.Cells.Item(4, 5).Value = 155578   # Recuperados
$ws.Cells.Item(4, 7).Value = 8        # Muertes hoy
$ws.Cells.Item(4, 8).Value = 3164     # Muertes

# Row 18: Canada
$ws.Cells.Item(18, 2).Value = 7474   # Casos totales
$ws.Cells.Item(18, 3).Value = 26     # Nuevos casos
$ws.Cells.Item(18, 4).Value = 1114   # Casos activos
$ws.Cells.Item(18, 5).Value = 6268   # Recuperados
$ws.Cells.Item(18, 7).Value = 3      # Muertes hoy
$ws.Cells.Item(18, 8).Value = 92     # Muertes

# Row 75: Tunez
$ws.Cells.Item(75, 5).Value = 350     # Recuperados
$ws.Cells.Item(75, 8).Value = 9       # Muertes

# --- Bump the "datos actualizados" footer timestamp ------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 04:20"
